# Apply "result step 3 and 4" changes:
#  - drop the sd_beta_log_bathymetry / sd_beta_mean_SST columns (Q:R)
#  - rename beta_log_bathymetry -> beta_mean_autumn_SST (col O)
#  - rename beta_mean_SST -> sd_beta_mean_autumn_SST (col P)
#  - refresh all model-selection figures (columns B:P) with the new run's numbers
#  - tighten the conditional-formatting thresholds for H:M

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Remove the last two columns (sd_beta_log_bathymetry, sd_beta_mean_SST)
$ws.Range("Q1:R5").EntireColumn.Delete()

# 2) Update the two remaining headers that changed meaning
$ws.Range("O1").Value = "beta_mean_autumn_SST"
$ws.Range("P1").Value = "sd_beta_mean_autumn_SST"

# 3) Refresh the data table (columns A:P, rows 2-5) with the new values
$rows = @(
    @("without spatial",  1.223, 109.8, 0.18, 0.51, 0.68, 0.5,  390, 580, 970, 387, 580, 967, 1, -1.26, 0.55),
    @("with spatial exp", 1.1,   33.9,  0.17, 0.52, 0.71, 0.5,  389, 578, 967, 396, 583, 979, 2, -1.34, 0.61),
    @("with sp shpere",   2.561, 49.2,  0.18, 0.51, 0.71, 0.51, 389, 579, 968, 396, 583, 979, 3, -1.32, 0.57),
    @("with sp gaussian", 6.603, 26.7,  0.18, 0.51, 0.68, 0.5,  390, 580, 970, 396, 583, 979, 4, -1.35, 0.62)
)

for ($r = 0; $r -lt $rows.Length; $r++) {
    $row = $rows[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r + 2, $c + 1).Value = $row[$c]
    }
}

# 4) Tighten the conditional-formatting expression thresholds
$ws.Range("H2:H5").FormatConditions.Item(2).Formula1 = "=H2<392"
$ws.Range("I2:I5").FormatConditions.Item(2).Formula1 = "=I2<581"
$ws.Range("J2:J5").FormatConditions.Item(2).Formula1 = "=J2<970"
$ws.Range("K2:K5").FormatConditions.Item(2).Formula1 = "=K2<390"
$ws.Range("L2:L5").FormatConditions.Item(2).Formula1 = "=L2<583"
$ws.Range("M2:M5").FormatConditions.Item(2).Formula1 = "=M2<970"
